$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated row in column A (the date column) and the row below it
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Clone the formatting of the last data row onto the new row (keeps the A-column
# date style + general text formatting for B:E identical to existing rows)
$ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 5)).Copy()
$ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 5)).PasteSpecial(-4122)

# New quote row for 2025-11-17 (date serial 45978)
$ws.Cells.Item($newRow, 1).Value = 45978
$ws.Cells.Item($newRow, 2).Value = "15,5317"
$ws.Cells.Item($newRow, 3).Value = "11,1166"
$ws.Cells.Item($newRow, 4).Value = "15,5317"
$ws.Cells.Item($newRow, 5).Value = "15,5317"
